# Apply the Alvearie FHIR IG metadata update (Version 5.0.0 -> 6.0.0, etc.)
# to the "Metadata" worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Old row 10 was "Contact" / "No display for ContactDetail";
# becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Old row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# it is replaced with the old "Description" row contents.
$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "IBM Health Data Connect standard values for organization type"

# Old row 12 ("Description" / long text) becomes "Purpose" with an empty value
$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = ""

# Old row 13 ("Purpose" / empty) becomes "Copyright" with an empty value
$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = ""

# Old row 14 ("Copyright" / empty) becomes "Immutable" / "BooleanType[null]"
$ws.Range("A14").Value = "Immutable"
$ws.Range("B14").Value = "BooleanType[null]"

# Old row 15 ("Immutable" / "BooleanType[null]") is now redundant and removed,
# shrinking the sheet from A1:B15 to A1:B14.
$ws.Rows.Item(15).Delete()
